# Finalize version to deliver reports to user.
# Update the counts in the SeenRx_CMT sheet (columns B:J, rows 2-42) with the
# final reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(52,25,11,5,12,30,23,1,14),
    @(6,3,0,0,0,3,7,0,2),
    @(0,0,0,0,0,0,0,0,0),
    @(6,3,0,0,0,1,3,0,2),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,1,0,0,0),
    @(0,0,0,0,0,1,4,0,0),
    @(9,3,1,0,2,1,7,0,6),
    @(4,0,0,0,0,0,0,0,1),
    @(0,0,0,0,0,1,1,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(2,2,0,0,2,0,1,0,0),
    @(3,1,1,0,0,0,4,0,3),
    @(0,0,0,0,0,0,1,0,1),
    @(3,3,1,0,0,1,0,0,0),
    @(3,3,1,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(7,3,4,1,3,3,0,1,2),
    @(1,0,2,1,3,1,0,0,1),
    @(3,3,2,0,0,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(1,0,0,0,0,1,0,1,0),
    @(2,0,0,0,0,0,0,0,1),
    @(7,4,3,0,0,10,7,0,2),
    @(4,2,1,0,0,5,2,0,1),
    @(2,0,2,0,0,2,4,0,1),
    @(0,1,0,0,0,2,0,0,0),
    @(1,1,0,0,0,1,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(20,9,2,4,7,12,2,0,2),
    @(7,3,0,2,4,2,2,0,1),
    @(4,3,0,2,0,4,0,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(5,0,1,0,3,4,0,0,0),
    @(4,3,1,0,0,2,0,0,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowNum = $startRow + $i
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($rowNum, $col + 2).Value = $rowValues[$col]
    }
}
